# Regenerate merged AHB files
# - rename the diff-header columns from the _old/_new naming convention
#   used by the previous FV release pair to the new FV2404/FV2410 pair
# - freeze the header row
# - turn the data range into a proper Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row ------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# "diff" header (column K) keeps its name.

# --- 2. Freeze the header row -------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn A1:U64 into an Excel Table --------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$lo.Name = "Table1"
